# Update examples to use the remote_access method with server_groups:
#  - column I (COMPLIANCE_GROUPS) is removed entirely
#  - the GROUPS header in H1 is renamed to SERVER_GROUPS
#  - row 3 (the SSH public-key row) grows much taller after the relayout
#  - the active selection moves to H1 (the renamed header cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the COMPLIANCE_GROUPS column outright - its only value ("Anssi") is
# not kept anywhere else in the sheet.
[void]$ws.Columns("I").Delete()

# Rename the GROUPS header to SERVER_GROUPS.
$ws.Range("H1").Value = "SERVER_GROUPS"

# Row 3 (server02 / SSH key row) is re-flowed to a much taller height.
$ws.Rows(3).RowHeight = 1303.7

# Selection ends up on the freshly renamed header cell.
[void]$ws.Range("H1").Select()

Write-Host "done"
